$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width: 44 -> 50 ---
$ws.Columns.Item(1).ColumnWidth = 50

# --- "Bad Drivers" table updates ---

# Row 3 (existing driver, values change)
$ws.Range("C3").Value = 704
$ws.Range("D3").Value = 83.8

# Row 4 (existing row, driver + values change)
$ws.Range("A4").Value = "Intel(R) Dual Band Wireless-AC 8260 - 20.70.25.2"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 97.90000000000001

# Insert a brand-new data row at position 5 (copy format from row 4, the last data row)
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.120.0.3"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 235
$ws.Range("D5").Value = 98.8

# Totals row (was row 5, now shifted to row 6)
$ws.Range("B6").Value = 7
$ws.Range("C6").Value = 940

# --- "Good Drivers" table updates ---
# Header now lives at row 13 (shifted down by 1 from the new row-5 insert above).
# Currently there are 6 data rows (14-19). We need 17 data rows (14-30),
# so insert 11 additional rows right after the header, inheriting the
# existing data-row formatting.
$ws.Range("A14:A24").EntireRow.Insert()

$ws.Range("A14").Value = "Intel(R) Dual Band Wireless-AC 8260 - 20.50.0.5"
$ws.Range("B14").Value = 323804
$ws.Range("D14").Value = 100

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B15").Value = 11128
$ws.Range("D15").Value = 100

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B16").Value = 486214
$ws.Range("D16").Value = 99.90000000000001
$ws.Range("E16").Value = "2024-11-10"

$ws.Range("A17").Value = "Intel(R) Dual Band Wireless-AC 8260 - 22.180.0.4"
$ws.Range("B17").Value = 10456
$ws.Range("D17").Value = 100
$ws.Range("E17").Value = "2022-10-17"

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B18").Value = 11140
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "2022-08-29"

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B19").Value = 14487
$ws.Range("D19").Value = 100
$ws.Range("E19").Value = "2022-05-23"

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B20").Value = 265400
$ws.Range("D20").Value = 99.90000000000001
$ws.Range("E20").Value = "2022-05-01"

$ws.Range("A21").Value = "Intel(R) Dual Band Wireless-AC 8260 - 22.80.1.1"
$ws.Range("B21").Value = 123675
$ws.Range("D21").Value = 100
$ws.Range("E21").Value = "2021-09-11"

$ws.Range("A22").Value = "Intel(R) Dual Band Wireless-AC 8260 - 20.70.27.1"
$ws.Range("B22").Value = 18967
$ws.Range("D22").Value = 100
$ws.Range("E22").Value = "2021-09-11"

$ws.Range("A23").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B23").Value = 79953
$ws.Range("D23").Value = 99.90000000000001
$ws.Range("E23").Value = "2021-08-18"

$ws.Range("A24").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B24").Value = 35355
$ws.Range("D24").Value = 100
$ws.Range("E24").Value = "2021-04-27"

$ws.Range("A25").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B25").Value = 65425
$ws.Range("D25").Value = 100
$ws.Range("E25").Value = "2020-08-05"

$ws.Range("A26").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B26").Value = 117653
$ws.Range("D26").Value = 100
$ws.Range("E26").Value = "2020-01-06"

$ws.Range("A27").Value = "Intel(R) Dual Band Wireless-AC 8260 - 20.70.16.4"
$ws.Range("B27").Value = 35023
$ws.Range("D27").Value = 100
$ws.Range("E27").Value = "2019-12-31"

$ws.Range("A28").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B28").Value = 56018
$ws.Range("D28").Value = 100
$ws.Range("E28").Value = "2019-12-14"

$ws.Range("A29").Value = "Intel(R) Dual Band Wireless-AC 8260 - 20.70.12.5"
$ws.Range("B29").Value = 197997
$ws.Range("D29").Value = 99.90000000000001
$ws.Range("E29").Value = "2019-08-25"

$ws.Range("A30").Value = "Intel(R) Dual Band Wireless-AC 8260 - 20.70.5.2"
$ws.Range("B30").Value = 160536
$ws.Range("D30").Value = 99.90000000000001
$ws.Range("E30").Value = "2018-11-25"
